# Update calculated price/profit values in Sheets per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 194.44444
$ws.Range("I4").Value = 168.75
$ws.Range("K4").Value = 168.75
$ws.Range("M4").Value = -54.75

$ws.Range("H88").Value = 8464.286
$ws.Range("J88").Value = 11670
$ws.Range("L88").Value = 11670
$ws.Range("N88").Value = -12482

$ws.Range("H91").Value = 8464.286
$ws.Range("J91").Value = 11670
$ws.Range("L91").Value = 11670
$ws.Range("N91").Value = -14478

$ws.Range("H98").Value = 3890.2
$ws.Range("I98").Value = 2086.5334
$ws.Range("J98").Value = 9301.200000000001
$ws.Range("K98").Value = 2086.5334
$ws.Range("L98").Value = 9301.200000000001
$ws.Range("M98").Value = -588.5333999999998
$ws.Range("N98").Value = -12297.2

$ws.Range("H111").Value = 950
$ws.Range("I111").Value = 900
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 2700
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = 367
$ws.Range("N111").Value = -9134

$ws.Range("H122").Value = 3890.2
$ws.Range("I122").Value = 2086.5334
$ws.Range("J122").Value = 9301.200000000001
$ws.Range("K122").Value = 6259.600199999999
$ws.Range("L122").Value = 27903.6
$ws.Range("M122").Value = -3809.600199999999
$ws.Range("N122").Value = -32803.60000000001

$ws.Range("H129").Value = 922.1070999999999
$ws.Range("J129").Value = 957.6226
$ws.Range("L129").Value = 2872.8678
$ws.Range("N129").Value = -12872.8678

$ws.Range("H132").Value = 6106.242
$ws.Range("I132").Value = 5767.643
$ws.Range("J132").Value = 8002.4
$ws.Range("K132").Value = 17302.929
$ws.Range("L132").Value = 24007.2
$ws.Range("M132").Value = -14772.929
$ws.Range("N132").Value = -29067.2

$ws.Range("H138").Value = 3129.861
$ws.Range("J138").Value = 3384.6
$ws.Range("L138").Value = 10153.8
$ws.Range("N138").Value = -20433.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5332.065
$ws.Range("I32").Value = 3325.9434
$ws.Range("K32").Value = 3325.9434
$ws.Range("M32").Value = -3038.9434

$ws.Range("H45").Value = 976.06665
$ws.Range("I45").Value = 817.2857
$ws.Range("J45").Value = 1115
$ws.Range("K45").Value = 817.2857
$ws.Range("L45").Value = 1115
$ws.Range("M45").Value = -440.2857
$ws.Range("N45").Value = -1869

$ws.Range("H61").Value = 2332.818
$ws.Range("I61").Value = 1875
$ws.Range("J61").Value = 2594.4285
$ws.Range("K61").Value = 1875
$ws.Range("L61").Value = 2594.4285
$ws.Range("M61").Value = -1663
$ws.Range("N61").Value = -3018.4285

$ws.Range("H122").Value = 2828.923
$ws.Range("I122").Value = 1655.7778
$ws.Range("J122").Value = 5468.5
$ws.Range("K122").Value = 4967.3334
$ws.Range("L122").Value = 16405.5
$ws.Range("M122").Value = -2517.3334
$ws.Range("N122").Value = -21305.5

$ws.Range("H136").Value = 2332.818
$ws.Range("I136").Value = 1875
$ws.Range("J136").Value = 2594.4285
$ws.Range("K136").Value = 5625
$ws.Range("L136").Value = 7783.2855
$ws.Range("M136").Value = -3075
$ws.Range("N136").Value = -12883.2855

$ws.Range("H137").Value = 39650
$ws.Range("J137").Value = 39650
$ws.Range("L137").Value = 39650
$ws.Range("N137").Value = -49850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 64608
$ws.Range("J59").Value = 64608
$ws.Range("L59").Value = 64608
$ws.Range("N59").Value = -66302

$ws.Range("H134").Value = 2894.7188
$ws.Range("I134").Value = 1279.2778
$ws.Range("K134").Value = 3837.8334
$ws.Range("M134").Value = -1302.8334

$ws.Range("H137").Value = 32500
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26320636
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 26320636
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 26320636
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -26321226

$ws.Range("H34").Value = 26320636
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 26320636
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 26320636
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -26321040

$ws.Range("H51").Value = 34073.4
$ws.Range("J51").Value = 34073.4
$ws.Range("L51").Value = 34073.4
$ws.Range("N51").Value = -35545.4

$ws.Range("H61").Value = 34073.4
$ws.Range("J61").Value = 34073.4
$ws.Range("L61").Value = 34073.4
$ws.Range("N61").Value = -34769.4

$ws.Range("H132").Value = 2553.6956
$ws.Range("I132").Value = 1557.6875
$ws.Range("K132").Value = 4673.0625
$ws.Range("M132").Value = -2143.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 4326.6665
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 4326.6665
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 12979.9995
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -13317.9995

$ws.Range("H30").Value = 4326.6665
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 4326.6665
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 12979.9995
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -13183.9995

$ws.Range("H113").Value = 787.0263
$ws.Range("I113").Value = 657.2083
$ws.Range("J113").Value = 1009.5714
$ws.Range("K113").Value = 1971.6249
$ws.Range("L113").Value = 3028.7142
$ws.Range("M113").Value = 198.3751
$ws.Range("N113").Value = -7368.7142

$ws.Range("H122").Value = 3680
$ws.Range("I122").Value = 1366.3334
$ws.Range("J122").Value = 3827.681
$ws.Range("K122").Value = 12297.0006
$ws.Range("L122").Value = 34449.129
$ws.Range("M122").Value = -9847.000599999999
$ws.Range("N122").Value = -39349.129

$ws.Range("H131").Value = 9438561
$ws.Range("I131").Value = 100040240
$ws.Range("J131").Value = 885.75
$ws.Range("K131").Value = 300120720
$ws.Range("L131").Value = 2657.25
$ws.Range("M131").Value = -300115680
$ws.Range("N131").Value = -12737.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 31039.285
$ws.Range("J46").Value = 31039.285
$ws.Range("L46").Value = 31039.285
$ws.Range("N46").Value = -31351.285

$ws.Range("H102").Value = 2710.6667
$ws.Range("I102").Value = 1879.8
$ws.Range("J102").Value = 3126.1
$ws.Range("K102").Value = 1879.8
$ws.Range("L102").Value = 3126.1
$ws.Range("M102").Value = -257.8
$ws.Range("N102").Value = -6370.1

$ws.Range("H113").Value = 1774.5
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340

$ws.Range("H132").Value = 3669.3462
$ws.Range("I132").Value = 2712
$ws.Range("J132").Value = 4371.4
$ws.Range("K132").Value = 8136
$ws.Range("L132").Value = 13114.2
$ws.Range("M132").Value = -5606
$ws.Range("N132").Value = -18174.2

$ws.Range("H136").Value = 23496.875
$ws.Range("J136").Value = 23496.875
$ws.Range("L136").Value = 70490.625
$ws.Range("N136").Value = -75590.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3421.121
$ws.Range("I7").Value = 1952.8823
$ws.Range("K7").Value = 1952.8823
$ws.Range("M7").Value = -1840.8823

$ws.Range("H122").Value = 4985.8423
$ws.Range("I122").Value = 3294
$ws.Range("J122").Value = 7312.125
$ws.Range("K122").Value = 9882
$ws.Range("L122").Value = 21936.375
$ws.Range("M122").Value = -7432
$ws.Range("N122").Value = -26836.375

$ws.Range("H126").Value = 3421.121
$ws.Range("I126").Value = 1952.8823
$ws.Range("K126").Value = 5858.6469
$ws.Range("M126").Value = -3388.6469

$ws.Range("H136").Value = 4859.294
$ws.Range("I136").Value = 2178.6667
$ws.Range("K136").Value = 6536.000100000001
$ws.Range("M136").Value = -3986.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 36230.727
$ws.Range("J123").Value = 36230.727
$ws.Range("L123").Value = 36230.727
$ws.Range("N123").Value = -46030.727

$ws.Range("I132").Value = 12910.777
$ws.Range("K132").Value = 38732.331
$ws.Range("M132").Value = -36202.331
